$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.044796083760568
$ws.Cells.Item(2, 4).Value = 1.045392648093783
$ws.Cells.Item(2, 5).Value = 0.992614727750844
$ws.Cells.Item(2, 6).Value = 1.043621970825725
$ws.Cells.Item(2, 9).Value = 1.043459900759641
$ws.Cells.Item(2, 10).Value = 1.04985925243749
$ws.Cells.Item(2, 11).Value = 1.048161015827129
$ws.Cells.Item(2, 12).Value = 0.9955398523335997
$ws.Cells.Item(2, 13).Value = 1.04639531853913
$ws.Cells.Item(2, 14).Value = 1.051350172996491

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.046148677058986
$ws.Cells.Item(3, 4).Value = 1.04642776934678
$ws.Cells.Item(3, 5).Value = 0.9936372048519299
$ws.Cells.Item(3, 6).Value = 1.045601348538653
$ws.Cells.Item(3, 9).Value = 1.043932582764505
$ws.Cells.Item(3, 10).Value = 1.050857253995145
$ws.Cells.Item(3, 11).Value = 1.049006985890745
$ws.Cells.Item(3, 12).Value = 0.9963617723202687
$ws.Cells.Item(3, 13).Value = 1.048182718305146
$ws.Cells.Item(3, 14).Value = 1.052349591830829

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.047021360145061
$ws.Cells.Item(4, 4).Value = 1.047095262747458
$ws.Cells.Item(4, 5).Value = 0.9942998659930998
$ws.Cells.Item(4, 6).Value = 1.046879323976043
$ws.Cells.Item(4, 9).Value = 1.044235748291527
$ws.Cells.Item(4, 10).Value = 1.051500048427794
$ws.Cells.Item(4, 11).Value = 1.04955148872771
$ws.Cells.Item(4, 12).Value = 0.9968940712668347
$ws.Cells.Item(4, 13).Value = 1.049336083775651
$ws.Cells.Item(4, 14).Value = 1.052993299105303

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04738763830506
$ws.Cells.Item(5, 4).Value = 1.047375332967939
$ws.Cells.Item(5, 5).Value = 0.994578699834602
$ws.Cells.Item(5, 6).Value = 1.047415928707279
$ws.Cells.Item(5, 9).Value = 1.044362558721632
$ws.Cells.Item(5, 10).Value = 1.051769573303285
$ws.Cells.Item(5, 11).Value = 1.049779710525772
$ws.Cells.Item(5, 12).Value = 0.9971179600053012
$ws.Cells.Item(5, 13).Value = 1.049820208180277
$ws.Cells.Item(5, 14).Value = 1.053263206737032

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.047449103227735
$ws.Cells.Item(6, 4).Value = 1.047422326275686
$ws.Cells.Item(6, 5).Value = 0.994625531979634
$ws.Cells.Item(6, 6).Value = 1.047505989035115
$ws.Cells.Item(6, 9).Value = 1.044383813314128
$ws.Cells.Item(6, 10).Value = 1.051814786533546
$ws.Cells.Item(6, 11).Value = 1.049817989851547
$ws.Cells.Item(6, 12).Value = 0.9971555583673455
$ws.Cells.Item(6, 13).Value = 1.049901451255848
$ws.Cells.Item(6, 14).Value = 1.053308484175265

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.047026256709513
$ws.Cells.Item(7, 4).Value = 1.047099007189104
$ws.Cells.Item(7, 5).Value = 0.994303590798249
$ws.Cells.Item(7, 6).Value = 1.046886496671275
$ws.Cells.Item(7, 9).Value = 1.044237445250024
$ws.Cells.Item(7, 10).Value = 1.051503652597071
$ws.Cells.Item(7, 11).Value = 1.049554540931139
$ws.Cells.Item(7, 12).Value = 0.9968970624462089
$ws.Cells.Item(7, 13).Value = 1.049342555598635
$ws.Cells.Item(7, 14).Value = 1.052996908392914

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.045253729451071
$ws.Cells.Item(8, 4).Value = 1.045742952142798
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.044291504150657
$ws.Cells.Item(8, 9).Value = 1.043620205506971
$ws.Cells.Item(8, 10).Value = 1.05019715306546
$ws.Cells.Item(8, 11).Value = 1.048447518831006
$ws.Cells.Item(8, 12).Value = 0.9958175282591056
$ws.Cells.Item(8, 13).Value = 1.047000051293578
$ws.Cells.Item(8, 14).Value = 1.051688553482111

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.042110481509635
$ws.Cells.Item(9, 4).Value = 1.043335509280773
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.039696448504831
$ws.Cells.Item(9, 9).Value = 1.042511759174496
$ws.Cells.Item(9, 10).Value = 1.047871784386871
$ws.Cells.Item(9, 11).Value = 1.046474341913634
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.042847014955785
$ws.Cells.Item(9, 14).Value = 1.049359882513274

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.040001071912495
$ws.Cells.Item(10, 4).Value = 1.041718101556942
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.036616939969973
$ws.Cells.Item(10, 9).Value = 1.041758564207379
$ws.Cells.Item(10, 10).Value = 1.04630550697576
$ws.Cells.Item(10, 11).Value = 1.045143393664752
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.040060335737761
$ws.Cells.Item(10, 14).Value = 1.047791380808585

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.039084244720728
$ws.Cells.Item(11, 4).Value = 1.04101470178102
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.035279402220576
$ws.Cells.Item(11, 9).Value = 1.041428990366504
$ws.Cells.Item(11, 10).Value = 1.045623384110342
$ws.Cells.Item(11, 11).Value = 1.044563314015446
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.038849179209102
$ws.Cells.Item(11, 14).Value = 1.047108289250456

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03874316630323
$ws.Cells.Item(12, 4).Value = 1.040752961251952
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.034781945110913
$ws.Cells.Item(12, 9).Value = 1.041306050855474
$ws.Cells.Item(12, 10).Value = 1.045369416056037
$ws.Cells.Item(12, 11).Value = 1.044347272418598
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.038398605106651
$ws.Cells.Item(12, 14).Value = 1.046853960532384

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.038816352856722
$ws.Cells.Item(13, 4).Value = 1.040809126726154
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.034888680479688
$ws.Cells.Item(13, 9).Value = 1.041332445459834
$ws.Cells.Item(13, 10).Value = 1.04542392024074
$ws.Cells.Item(13, 11).Value = 1.044393640182971
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.038495286623585
$ws.Cells.Item(13, 14).Value = 1.046908542119281

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.039056061908829
$ws.Cells.Item(14, 4).Value = 1.040993075791408
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.035238295361029
$ws.Cells.Item(14, 9).Value = 1.041418838807691
$ws.Cells.Item(14, 10).Value = 1.045602403268812
$ws.Cells.Item(14, 11).Value = 1.044545467708025
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.038811948939025
$ws.Cells.Item(14, 14).Value = 1.047087278613726

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.039203684264424
$ws.Cells.Item(15, 4).Value = 1.041106350818299
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.035453619699181
$ws.Cells.Item(15, 9).Value = 1.041471999398385
$ws.Cells.Item(15, 10).Value = 1.04571229313592
$ws.Cells.Item(15, 11).Value = 1.044638937338544
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.039006962139108
$ws.Cells.Item(15, 14).Value = 1.047197324537049

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.040061845296479
$ws.Cells.Item(16, 4).Value = 1.041764718903892
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.036705619937105
$ws.Cells.Item(16, 9).Value = 1.041780364133503
$ws.Cells.Item(16, 10).Value = 1.046350693972683
$ws.Cells.Item(16, 11).Value = 1.045181811546144
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.040140619597649
$ws.Cells.Item(16, 14).Value = 1.047836631976226

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.040599218392741
$ws.Cells.Item(17, 4).Value = 1.042176873152797
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.037489856417636
$ws.Cells.Item(17, 9).Value = 1.041972869923656
$ws.Cells.Item(17, 10).Value = 1.046750091522209
$ws.Cells.Item(17, 11).Value = 1.045521327343694
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.040850512766212
$ws.Cells.Item(17, 14).Value = 1.048236596716084

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.040912327993958
$ws.Cells.Item(18, 4).Value = 1.042416981872518
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.037946895106214
$ws.Cells.Item(18, 9).Value = 1.042084824162318
$ws.Cells.Item(18, 10).Value = 1.046982676295807
$ws.Cells.Item(18, 11).Value = 1.045718997915534
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.041264148012572
$ws.Cells.Item(18, 14).Value = 1.048469511786739

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.041019034467718
$ws.Cells.Item(19, 4).Value = 1.042498803136295
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.038102667445217
$ws.Cells.Item(19, 9).Value = 1.042122941648363
$ws.Cells.Item(19, 10).Value = 1.047061918012432
$ws.Cells.Item(19, 11).Value = 1.045786337091842
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.041405114016835
$ws.Cells.Item(19, 14).Value = 1.04854886603569

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.040541597667366
$ws.Cells.Item(20, 4).Value = 1.04213268335826
$ws.Cells.Item(20, 5).Value = 0.9894336180355766
$ws.Cells.Item(20, 6).Value = 1.037405756089228
$ws.Cells.Item(20, 9).Value = 1.041952250160689
$ws.Cells.Item(20, 10).Value = 1.046707279012205
$ws.Cells.Item(20, 11).Value = 1.045484938140253
$ws.Cells.Item(20, 12).Value = 0.9929783193490043
$ws.Cells.Item(20, 13).Value = 1.040774392920999
$ws.Cells.Item(20, 14).Value = 1.048193723407405

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.038985488233841
$ws.Cells.Item(21, 4).Value = 1.040938920320399
$ws.Cells.Item(21, 5).Value = 0.9882828385668255
$ws.Cells.Item(21, 6).Value = 1.035135360224155
$ws.Cells.Item(21, 9).Value = 1.041393412545535
$ws.Cells.Item(21, 10).Value = 1.045549861026698
$ws.Cells.Item(21, 11).Value = 1.044500774199119
$ws.Cells.Item(21, 12).Value = 0.9920501090198107
$ws.Cells.Item(21, 13).Value = 1.03871871916941
$ws.Cells.Item(21, 14).Value = 1.047034661755601

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.038004040062652
$ws.Cells.Item(22, 4).Value = 1.040185650832578
$ws.Cells.Item(22, 5).Value = 0.9875604150241496
$ws.Cells.Item(22, 6).Value = 1.033704181487999
$ws.Cells.Item(22, 9).Value = 1.041039031655232
$ws.Cells.Item(22, 10).Value = 1.044818686282884
$ws.Cells.Item(22, 11).Value = 1.043878665297437
$ws.Cells.Item(22, 12).Value = 0.991467000034148
$ws.Cells.Item(22, 13).Value = 1.037422195043636
$ws.Cells.Item(22, 14).Value = 1.046302448659784

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038524618210234
$ws.Cells.Item(23, 4).Value = 1.040585232199644
$ws.Cells.Item(23, 5).Value = 0.9879432794636459
$ws.Cells.Item(23, 6).Value = 1.034463233503578
$ws.Cells.Item(23, 9).Value = 1.041227183370403
$ws.Cells.Item(23, 10).Value = 1.045206626911427
$ws.Cells.Item(23, 11).Value = 1.044208774869276
$ws.Cells.Item(23, 12).Value = 0.9917760702887607
$ws.Cells.Item(23, 13).Value = 1.038109896622432
$ws.Cells.Item(23, 14).Value = 1.046690940208517

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.04056763501581
$ws.Cells.Item(24, 4).Value = 1.042152651732438
$ws.Cells.Item(24, 5).Value = 0.9894529299347241
$ws.Cells.Item(24, 6).Value = 1.037443758622109
$ws.Cells.Item(24, 9).Value = 1.041961568367194
$ws.Cells.Item(24, 10).Value = 1.046726625307731
$ws.Cells.Item(24, 11).Value = 1.045501381977372
$ws.Cells.Item(24, 12).Value = 0.9929938892766438
$ws.Cells.Item(24, 13).Value = 1.040808789541447
$ws.Cells.Item(24, 14).Value = 1.04821309717689

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042925495492864
$ws.Cells.Item(25, 4).Value = 1.043960054247836
$ws.Cells.Item(25, 5).Value = 0.9912096547607046
$ws.Cells.Item(25, 6).Value = 1.040887142679539
$ws.Cells.Item(25, 9).Value = 1.042800808732039
$ws.Cells.Item(25, 10).Value = 1.048475739047921
$ws.Cells.Item(25, 11).Value = 1.046987157705932
$ws.Cells.Item(25, 12).Value = 0.9944092447426411
$ws.Cells.Item(25, 13).Value = 1.043923769877287
$ws.Cells.Item(25, 14).Value = 1.049964694859217

Write-Host "Applied vm_pu updates"